$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.286.65"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.030.45"
$ws.Range("E3").Value = "  -0.81%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'243.99"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").Value = "'0.655"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'52.54"
$ws.Range("E8").Value = "  -7.93%  "
$ws.Range("D9").Value = "'60.96"
$ws.Range("E9").Value = "  +4.06%  "
$ws.Range("D10").Value = "'0.357"
$ws.Range("E10").Value = "  -2.90%  "
$ws.Range("D11").Value = "'0.0735"
$ws.Range("E11").Value = "  -4.96%  "
$ws.Range("D12").Value = "'0.105"
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").Value = "'0.929"
$ws.Range("E13").Value = "  +6.98%  "
$ws.Range("D14").Value = "'14.27"
$ws.Range("E14").Value = "  -5.69%  "
$ws.Range("D15").Value = "2.330.38"
$ws.Range("E15").Value = "  -0.67%  "
$ws.Range("D16").Value = "'5.30"
$ws.Range("E16").Value = "  -5.08%  "
$ws.Range("D17").Value = "2.039.06"
$ws.Range("E17").Value = "  -0.81%  "
$ws.Range("D18").Value = "36.188.82"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "'16.74"
$ws.Range("E19").Value = "  -6.10%  "
$ws.Range("D20").Value = "'70.67"
$ws.Range("E20").Value = "  -3.51%  "
$ws.Range("D21").Value = "0.0₃0840"
$ws.Range("E21").Value = "  -4.82%  "
$ws.Range("D22").Value = "'235.44"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Value = "'5.10"
$ws.Range("E23").Value = "  -4.90%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  -3.79%  "
$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "'163.15"
$ws.Range("E27").Value = "  -2.62%  "
$ws.Range("E28").Value = "  -11.33%  "
$ws.Range("D29").Value = "'19.60"
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D31").Value = "'1.16"
$ws.Range("E31").Value = "  +5.10%  "
$ws.Range("D32").Value = "'4.88"
$ws.Range("E32").Value = "  -10.07%  "
$ws.Range("D33").Value = "'0.0584"
$ws.Range("E33").Value = "  -3.97%  "
$ws.Range("D34").Value = "'4.31"
$ws.Range("E34").Value = "  -8.59%  "
$ws.Range("E35").Value = "  +0.04%  "
$ws.Range("D36").Value = "'0.0845"
$ws.Range("E36").Value = "  +4.63%  "
$ws.Range("E37").Value = "  -0.93%  "
$ws.Range("E38").Value = "  -5.10%  "
$ws.Range("D39").Value = "'4.88"
$ws.Range("E39").Value = "  -3.86%  "
$ws.Range("E40").Value = "  -7.62%  "
$ws.Range("D41").Value = "'2.87"
$ws.Range("E41").Value = "  -4.73%  "
$ws.Range("D42").Value = "'0.0210"
$ws.Range("E42").Value = "  -5.35%  "
$ws.Range("E43").Value = "  -5.28%  "
$ws.Range("D44").Value = "'91.69"
$ws.Range("E44").Value = "  -4.59%  "
$ws.Range("D45").Value = "'0.0883"
$ws.Range("E45").Value = "  -5.82%  "
$ws.Range("D46").Value = "1.371.86"
$ws.Range("E46").Value = "  +5.33%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'7.37"
$ws.Range("E47").Value = "  +9.57%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'15.53"
$ws.Range("E48").Value = "  -7.49%  "
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "2.218.20"
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("D51").Value = "'2.22"
$ws.Range("E51").Value = "  -4.96%  "

# Reset style on cells where we used the text-prefix trick, so no stray
# number-format/quotePrefix style survives on the cell (matches original formatting)
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D51").Style = "Normal"
